# Regenerate the "K" column (G) of the save_data sheet.
# The column used to hold a raw strike count (Strike#); it is being
# recomputed here and rewritten with the new K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..39 (one row per recorded start), in order.
$kVals = @(12, 8, 5, 3, 7, 3, 5, 6, 6, 6, 7, 6, 5, 5, 4, 5, 12, 4, 2, 5, 4, 8, 5, 4, 7, 8, 5, 6, 4, 5, 7, 5, 0, 4, 1, 5, 3, 0)

$firstRow = 2
$lastRow = 39

# Build a 2D (rows x 1 column) array for a single bulk write into G2:G39.
$rowCount = $lastRow - $firstRow + 1
$data = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i, 0] = $kVals[$i]
}

$rng = $ws.Range("G$firstRow`:G$lastRow")
$rng.Value = $data
